# "opción de medrar a vip" — mark the two "Modificar datos personales" (FER)
# subtasks as done ("Listo (sin testeo)", green) instead of
# "Lógica hecha, falta GUI y conectarlas" (amber/orange), and move the
# sheet's visible selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

foreach ($addr in @("E8", "E13")) {
    $cell = $ws.Range($addr)
    $cell.Value = "Listo (sin testeo)"
    $cell.Interior.Color = 5296274   # RGB(146,208,80) == fill FF92D050 (same as E4/E19/E21/E22)
}

# Move the active selection / view to match the edited sheet state.
$null = $ws.Range("G15").Select()
